# Applies the betexplorer scraper update for
# 2023/israel_ligat-ha-al_2023-2024.xlsx
#
# 1) Swaps the F:V match-data for four row pairs (the scraper re-ordered
#    matches that share the same match date): 16<->17, 18<->19, 24<->25,
#    31<->32.
# 2) Appends one brand-new match row (row 36 / Indice 35).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-MatchRows($rowA, $rowB) {
    # Columns F..V (6..22) hold the per-match data; A..E (index/pais/
    # torneio/temporada/data_partida) are identical between the two rows
    # in every pair, so only F:V needs to move.
    for ($col = 6; $col -le 22; $col++) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)

        $valueA = $cellA.Value2
        $valueB = $cellB.Value2

        $cellA.Value = $valueB
        $cellB.Value = $valueA
    }
}

Swap-MatchRows 16 17
Swap-MatchRows 18 19
Swap-MatchRows 24 25
Swap-MatchRows 31 32

# New row 36 appended at the bottom of the table.
# Column A ("Indice") uses the bold/bordered/centered style and column E
# ("data_partida") uses the datetime number format throughout the sheet,
# so copy those formats down from the previous row (35) before writing
# the new values, rather than leaving the new cells with the default
# style (which would mint new style entries in styles.xml).
$r = 36
$ws.Cells.Item(35, 1).Copy() | Out-Null
$ws.Cells.Item($r, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(35, 5).Copy() | Out-Null
$ws.Cells.Item($r, 5).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($r, 1).Value = 35
$ws.Cells.Item($r, 2).Value = "israel"
$ws.Cells.Item($r, 3).Value = "ligat-ha-al"
$ws.Cells.Item($r, 4).Value = "2023-2024"
$ws.Cells.Item($r, 5).Value = 45255.79166666666
$ws.Cells.Item($r, 6).Value = "Maccabi Haifa"
$ws.Cells.Item($r, 7).Value = 2
$ws.Cells.Item($r, 8).Value = "Hapoel Petah Tikva"
$ws.Cells.Item($r, 9).Value = 1
$ws.Cells.Item($r, 10).Value = 1.15
$ws.Cells.Item($r, 11).Value = "25/11/2023 18:29"
$ws.Cells.Item($r, 12).Value = 1.16
$ws.Cells.Item($r, 13).Value = "25/11/2023 18:43"
$ws.Cells.Item($r, 14).Value = 8.56
$ws.Cells.Item($r, 15).Value = "25/11/2023 18:29"
$ws.Cells.Item($r, 16).Value = 7.56
$ws.Cells.Item($r, 17).Value = "25/11/2023 18:45"
$ws.Cells.Item($r, 18).Value = 16.61
$ws.Cells.Item($r, 19).Value = "25/11/2023 18:29"
$ws.Cells.Item($r, 20).Value = 14.19
$ws.Cells.Item($r, 21).Value = "25/11/2023 18:45"
$ws.Cells.Item($r, 22).Value = "https://www.betexplorer.com/football/israel/ligat-ha-al/maccabi-haifa-hapoel-petah-tikva/hxrcfAm6/"
